$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of data (rows 5-11) following the existing 4 rows
$ws.Range("A5").Value = "TechGig Webinar <expertspeak@techgig.com>"

$ws.Range("A6").Value = "Mailtrack Notification <notification@mailtrack.io>"

$ws.Range("A7").Value = "116_Deepak Kumar <deepakkumar737373@gmail.com>"
$ws.Range("B7").Value = "non dust eraser"

$ws.Range("A8").Value = "116_Deepak Kumar <deepakkumar737373@gmail.com>"
$ws.Range("B8").Value = "non dust eraser"

$ws.Range("A9").Value = "116_Deepak Kumar <deepakkumar737373@gmail.com>"
$ws.Range("B9").Value = "non dust eraser"

$ws.Range("A10").Value = "kundan Prakash Jha <kundan.ext123@gmail.com>"
$ws.Range("B10").Value = "Operating system"

$ws.Range("A11").Value = "Mailtrack Notification <notification@mailtrack.io>"

# Column widths to match target layout (values chosen so the runtime's
# pixel-quantized ColumnWidth -> XML width conversion lands on the closest
# achievable width to the target 57.77734375 / 49.44140625)
$ws.Columns.Item(1).ColumnWidth = 56.92
$ws.Columns.Item(2).ColumnWidth = 48.6

# Selection moved to A6, matching target sheetView
$null = $ws.Range("A6").Select()
